# Add the new game record (#13, played 2025-10-14) as row 14 on the
# Game_Record sheet. Game_Record!A2:A13 already hold a "ROW()-1" shared
# formula (si="0"); row 8's formula anchors a second shared group
# (si="1") that Excel will auto-extend to A14 when the formula is
# filled/entered the same way as its neighbours.

$wb = $excel.ActiveWorkbook
$gameRecord = $wb.Worksheets.Item("Game_Record")
$statSheet  = $wb.Worksheets.Item("Stat_Sheet")

$gameRecord.Activate()

$gameRecord.Range("A14").Formula = "=ROW()-1"
$gameRecord.Range("B14").Value = 45944
$gameRecord.Range("C14").Value = "SiderFace"
$gameRecord.Range("D14").Value = "SimpleJack"
$gameRecord.Range("E14").Value = "DrSystomatix"
$gameRecord.Range("F14").Value = "Player1"

# Match the author's recorded selections after the edit: Game_Record
# is left with A15 selected (just past the new row) and Stat_Sheet
# keeps K4 selected, with Game_Record as the active tab.
$statSheet.Activate()
$statSheet.Range("K4").Select()

$gameRecord.Activate()
$gameRecord.Range("A15").Select()

$wb.Save()
